# The post stored in row 744 ("「お金では買えないよ」...") was removed from the
# source data. Delete that entire row and shift all subsequent rows up by one,
# which matches the diff (dimension shrinks from A1:C801 to A1:C800 and every
# row after the old 744 is renumbered down by one).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("744").Delete()
